$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time formatting from the last existing data row (24) down to the
# two new rows so the same style indices (date + time number formats) get reused.
$ws.Range("A24:C24").Copy()
$ws.Range("A25:C26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Row 25
$ws.Range("A25").Value = 42115
$ws.Range("B25").Value = 0.54166666666666663
$ws.Range("C25").Value = 0.64583333333333337
$ws.Range("D25").Value = "Juliano, Tommy"
$ws.Range("E25").Value = "Einbauen einer Schachbrett-View um das Spiel mehrspielerfähig zu machne"

# Row 26
$ws.Range("A26").Value = 42120
$ws.Range("B26").Value = 0.45833333333333331
$ws.Range("C26").Value = 0.625
$ws.Range("D26").Value = "Juliano"
$ws.Range("E26").Value = "Schachbrett-View-Anpassungen"
$ws.Range("F26").Value = " "

# Update sheet view: scroll so row 7 is at the top-left, and select D26
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D26").Select()
